# Auto-generated edit script: updates H..N numeric columns across several
# sheets to reflect a scheduled market-data refresh (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 413.52777
$ws.Range("I19").Value = 412.44446
$ws.Range("J19").Value = 414.6111
$ws.Range("K19").Value = 412.44446
$ws.Range("L19").Value = 414.6111
$ws.Range("M19").Value = -237.44446
$ws.Range("N19").Value = -764.6111000000001
$ws.Range("H41").Value = 570.1429000000001
$ws.Range("I41").Value = 269.7143
$ws.Range("J41").Value = 870.5714
$ws.Range("K41").Value = 269.7143
$ws.Range("L41").Value = 870.5714
$ws.Range("M41").Value = 170.2857
$ws.Range("N41").Value = -1750.5714
$ws.Range("H43").Value = 2300.2
$ws.Range("I43").Value = 2125.25
$ws.Range("J43").Value = 3000
$ws.Range("K43").Value = 2125.25
$ws.Range("L43").Value = 3000
$ws.Range("M43").Value = -2056.25
$ws.Range("N43").Value = -3138
$ws.Range("H51").Value = 2250
$ws.Range("I51").Value = 2000
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 2000
$ws.Range("L51").Value = 3000
$ws.Range("M51").Value = -1516
$ws.Range("N51").Value = -3968
$ws.Range("H98").Value = 1222.8077
$ws.Range("I98").Value = 1009.35
$ws.Range("J98").Value = 1934.3334
$ws.Range("K98").Value = 1009.35
$ws.Range("L98").Value = 1934.3334
$ws.Range("M98").Value = 488.65
$ws.Range("N98").Value = -4930.3334
$ws.Range("H122").Value = 1222.8077
$ws.Range("I122").Value = 1009.35
$ws.Range("J122").Value = 1934.3334
$ws.Range("K122").Value = 3028.05
$ws.Range("L122").Value = 5803.0002
$ws.Range("M122").Value = -578.0500000000002
$ws.Range("N122").Value = -10703.0002
$ws.Range("H132").Value = 2142.087
$ws.Range("I132").Value = 1783.5
$ws.Range("K132").Value = 5350.5
$ws.Range("M132").Value = -2820.5
$ws.Range("H135").Value = 1922.3226
$ws.Range("I135").Value = 1688.5686
$ws.Range("J135").Value = 3006.0908
$ws.Range("K135").Value = 15197.1174
$ws.Range("L135").Value = 27054.8172
$ws.Range("M135").Value = -12662.1174
$ws.Range("N135").Value = -32124.8172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7162.63
$ws.Range("I32").Value = 5600.1953
$ws.Range("J32").Value = 17618.924
$ws.Range("K32").Value = 5600.1953
$ws.Range("L32").Value = 17618.924
$ws.Range("M32").Value = -5313.1953
$ws.Range("N32").Value = -18192.924
$ws.Range("H61").Value = 204422.28
$ws.Range("I61").Value = 5033.8887
$ws.Range("J61").Value = 438486.9
$ws.Range("K61").Value = 5033.8887
$ws.Range("L61").Value = 438486.9
$ws.Range("M61").Value = -4821.8887
$ws.Range("N61").Value = -438910.9
$ws.Range("H97").Value = 1184.0454
$ws.Range("I97").Value = 1123.579
$ws.Range("K97").Value = 1123.579
$ws.Range("M97").Value = -627.579
$ws.Range("H132").Value = 2276324.8
$ws.Range("I132").Value = 2608.9312
$ws.Range("K132").Value = 7826.7936
$ws.Range("M132").Value = -5296.7936
$ws.Range("H136").Value = 204422.28
$ws.Range("I136").Value = 5033.8887
$ws.Range("J136").Value = 438486.9
$ws.Range("K136").Value = 15101.6661
$ws.Range("L136").Value = 1315460.7
$ws.Range("M136").Value = -12551.6661
$ws.Range("N136").Value = -1320560.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 9787.968000000001
$ws.Range("I105").Value = 14625
$ws.Range("K105").Value = 14625
$ws.Range("M105").Value = -12878
$ws.Range("H134").Value = 51970.91
$ws.Range("I134").Value = 7234
$ws.Range("J134").Value = 335304.66
$ws.Range("K134").Value = 21702
$ws.Range("L134").Value = 1005913.98
$ws.Range("M134").Value = -19167
$ws.Range("N134").Value = -1010983.98

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 333341340
$ws.Range("H27").Value = 333341340
$ws.Range("H94").Value = 6573.737
$ws.Range("J94").Value = 6662.091
$ws.Range("L94").Value = 6662.091
$ws.Range("N94").Value = -7564.091
$ws.Range("H122").Value = 927249
$ws.Range("I122").Value = 1323919.9
$ws.Range("J122").Value = 1683.6666
$ws.Range("K122").Value = 3971759.7
$ws.Range("L122").Value = 5050.9998
$ws.Range("M122").Value = -3969309.7
$ws.Range("N122").Value = -9950.9998
$ws.Range("H134").Value = 280798.62
$ws.Range("I134").Value = 3042.4333
$ws.Range("K134").Value = 9127.2999
$ws.Range("M134").Value = -6592.2999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1215.5526
$ws.Range("J5").Value = 1685.625
$ws.Range("L5").Value = 5056.875
$ws.Range("N5").Value = -5280.875
$ws.Range("H86").Value = 894.2857
$ws.Range("I86").Value = 894.2857
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2682.8571
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1496.8571
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 894.2857
$ws.Range("I89").Value = 894.2857
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 8048.571300000001
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -2120.571300000001
$ws.Range("N89").ClearContents()
$ws.Range("H97").Value = 2857403.2
$ws.Range("I97").Value = 4000239.8
$ws.Range("J97").Value = 312.2
$ws.Range("K97").Value = 12000719.4
$ws.Range("L97").Value = 936.5999999999999
$ws.Range("M97").Value = -12000223.4
$ws.Range("N97").Value = -1928.6
$ws.Range("H122").Value = 3866.606
$ws.Range("I122").Value = 380.2381
$ws.Range("J122").Value = 9967.75
$ws.Range("K122").Value = 3422.1429
$ws.Range("L122").Value = 89709.75
$ws.Range("M122").Value = -972.1428999999998
$ws.Range("N122").Value = -94609.75
$ws.Range("H132").Value = 2979.375
$ws.Range("I132").Value = 666.6667
$ws.Range("J132").Value = 4367
$ws.Range("K132").Value = 6000.0003
$ws.Range("L132").Value = 39303
$ws.Range("M132").Value = -3470.0003
$ws.Range("N132").Value = -44363
$ws.Range("H135").Value = 1215.5526
$ws.Range("J135").Value = 1685.625
$ws.Range("L135").Value = 15170.625
$ws.Range("N135").Value = -20240.625
$ws.Range("H140").Value = 2209.524
$ws.Range("I140").Value = 2209.524
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 6628.572
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -1448.572
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 505500
$ws.Range("I35").Value = 1000000
$ws.Range("J35").Value = 11000
$ws.Range("K35").Value = 1000000
$ws.Range("L35").Value = 11000
$ws.Range("M35").Value = -999702
$ws.Range("N35").Value = -11596
$ws.Range("H70").Value = 5427.76
$ws.Range("I70").Value = 5583
$ws.Range("J70").Value = 5028.5713
$ws.Range("K70").Value = 5583
$ws.Range("L70").Value = 5028.5713
$ws.Range("M70").Value = -5313
$ws.Range("N70").Value = -5568.5713
$ws.Range("H73").Value = 5427.76
$ws.Range("I73").Value = 5583
$ws.Range("J73").Value = 5028.5713
$ws.Range("K73").Value = 5583
$ws.Range("L73").Value = 5028.5713
$ws.Range("M73").Value = -4647
$ws.Range("N73").Value = -6900.5713
$ws.Range("H113").Value = 38462932
$ws.Range("I113").Value = 83334350
$ws.Range("J113").Value = 1715.0714
$ws.Range("K113").Value = 83334350
$ws.Range("L113").Value = 1715.0714
$ws.Range("M113").Value = -83332180
$ws.Range("N113").Value = -6055.0714
$ws.Range("H132").Value = 2183.1025
$ws.Range("I132").Value = 2018.1923
$ws.Range("J132").Value = 2512.923
$ws.Range("K132").Value = 6054.5769
$ws.Range("L132").Value = 7538.768999999999
$ws.Range("M132").Value = -3524.5769
$ws.Range("N132").Value = -12598.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 2800
$ws.Range("J38").Value = 2800
$ws.Range("L38").Value = 2800
$ws.Range("N38").Value = -3620
$ws.Range("H55").Value = 176.78572
$ws.Range("I55").Value = 167.5
$ws.Range("J55").Value = 200
$ws.Range("K55").Value = 167.5
$ws.Range("L55").Value = 200
$ws.Range("M55").Value = 5.5
$ws.Range("N55").Value = -546
$ws.Range("H136").Value = 7394.617
$ws.Range("I136").Value = 6652.4614
$ws.Range("J136").Value = 8313.477000000001
$ws.Range("K136").Value = 19957.3842
$ws.Range("L136").Value = 24940.431
$ws.Range("M136").Value = -17407.3842
$ws.Range("N136").Value = -30040.431

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1233.7046
$ws.Range("I113").Value = 837.7714
$ws.Range("J113").Value = 2773.4443
$ws.Range("K113").Value = 2513.3142
$ws.Range("L113").Value = 8320.332900000001
$ws.Range("M113").Value = -343.3141999999998
$ws.Range("N113").Value = -12660.3329
$ws.Range("H132").Value = 1788.8235
$ws.Range("I132").Value = 1200.8462
$ws.Range("J132").Value = 3699.75
$ws.Range("K132").Value = 3602.5386
$ws.Range("L132").Value = 11099.25
$ws.Range("M132").Value = -1072.5386
$ws.Range("N132").Value = -16159.25
$ws.Range("H136").Value = 2648.628
$ws.Range("I136").Value = 2867.423
$ws.Range("J136").Value = 2314
$ws.Range("K136").Value = 8602.269
$ws.Range("L136").Value = 6942
$ws.Range("M136").Value = -6052.269
$ws.Range("N136").Value = -12042
